# Implemented and debugged CheckManageFile(), CheckTranslocFile() and started to work on ReadManageFile()
#
# Updates the "Description" and "TranslocationFile" sheets of the translocation
# file template: the Source/Target patch columns are renamed from
# "SourcePatchID"/"TargetPatchID" to the simpler "Source"/"Target", and the
# VALUES column in the Description sheet is expanded to note that a
# semicolon-separated list of integers is now also accepted.

$wb = $excel.ActiveWorkbook

$wsDescription = $wb.Worksheets.Item("Description")
$wsTransloc    = $wb.Worksheets.Item("TranslocationFile")

# --- Description sheet -----------------------------------------------------
# Row 4: Source patch id / location
$wsDescription.Range("A4").Value = "Source"
$wsDescription.Range("C4").Value = "Integer value or semicolon seperated integer values (cell-based)"

# Row 5: Target patch id / location
$wsDescription.Range("A5").Value = "Target"
$wsDescription.Range("C5").Value = "Integer value or semicolon seperated integer values (cell-based)"

# --- TranslocationFile sheet -------------------------------------------------
$wsTransloc.Range("C1").Value = "Source"
$wsTransloc.Range("D1").Value = "Target"

# --- Selection / active sheet bookkeeping -----------------------------------
# Description becomes the active sheet, with the bottom-right pane focused on A5.
$wsDescription.Activate()
$wsDescription.Range("A5").Select()

# TranslocationFile keeps a selection on its now-renamed "Target" header cell.
$wsTransloc.Range("D1").Select()

# Re-activate Description so it is the workbook's active/visible tab on save.
$wsDescription.Activate()
